$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.768452167510986
$ws.Range("B1").Value = 6.651125431060791
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 3.796440839767456
$ws.Range("E1").Value = 1.767792820930481
